# Standardize the weeks (add a new week column) and update the journal stubs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N (14) for the new week (4/21) -- this shifts the
# old N/O columns (4/24 and the "+14" formula week) out to O/P, carrying
# their styles/formulas/values along automatically.
$ws.Columns.Item(14).Insert()

# New week's date header.
$ws.Range("N2").Value = 44307

# New week's attendance values for each attendee row.
$ws.Range("N3").Value = 1
$ws.Range("N4").Value = 1
$ws.Range("N5").Value = 1
$ws.Range("N6").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("N8").Value = 1

# Fill in the two previously-unused weeks (L, M) that now have attendance
# data recorded for them as well.
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1

$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1

$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1

$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1

$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1

$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0

# Extend the weekly totals to include the new week column N.
$ws.Range("B3").Formula = "=SUM(C3:N3)"
$ws.Range("B4:B8").Formula = "=SUM(C4:N4)"

# Rename the "Maximum" journal stub to be more descriptive.
$ws.Range("A9").Value = "Maximum Possible (Excluding optional meetings)"

# Widen column A to fit the new, longer label and mark it best-fit.
$ws.Columns.Item(1).ColumnWidth = 43.14

# Clean up the view: drop the stale scrolled top-left cell and move the
# active selection to D10.
$ws.Range("D10").Select()
